$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value2 = 4635434752
$ws.Range("B3").Value2 = @'
Zora
'@
$ws.Range("C3").Value2 = @'
https://dierenasielgenk.be/wp-content/uploads/2021/07/IMG_6238-1200x800.jpg
'@
$ws.Range("D3").Value2 = @'
Via een inbeslagname is Zora in ons asiel binnengekomen en uiteindelijk ook toegewezen.
Ze heeft tijd nodig gehad om de mensen van het asiel te leren kennen, vanaf het moment dat het ijs was gebroken, is Zora een erg aanhankelijke hond voor de mensen die ze kent.
Werken wil Zora wel en spelen is ook zeker een hobby. Je hebt haar volledige aandacht wanneer je lekkere vleesjes hebt en dan merk je ook dat ze de basiscommando’s van thuis uit kent. Ze kent de commando’s ‘zit, ‘liggen’, ‘poot’ en ‘blijf’ en gezien haar leergierigheid, weten we zeker dat ze véél meer commando’s kan leren.
Zora wil ook graag samen spelen en als ze wil dat je iets gooit komt ze het speeltje op je schoot leggen.
Zora wordt wel het best met rust gelaten tijdens het eten of wanneer ze lekker rustig ligt te knabbelen op een kauwbotje. Ze heeft moeite met vreemde mensen en andere dieren dus hier moet de nieuwe eigenaar zeker rekening mee houden. We zoeken een thuis met mensen die het Mechelse herder ras kennen. Daarnaast een huis met een goed afgesloten tuin zonder kinderen of andere dieren.

'@
$ws.Range("E3").Value2 = @'
Vicky
'@
$ws.Range("F3").Value2 = @'
https://dierenasielgenk.be/wp-content/uploads/2021/04/vickyyy-1-1200x800.jpg
'@
$ws.Range("G3").Value2 = @'
Er zijn zo van die honden die bijna horen tot het ‘vaste meubilair’ van het asiel. Zonde, want we begrijpen vaak zelf niet goed hoe dat komt.
Een van die honden is Mechelaar Vicky.
Deze 7-jarige dame zit al bijna 3 jaar in ons asiel en ze zit samen met Molly in onze top drie van langzitters. Haar blinkende ogen verraden weinig over haar trieste verleden. Vicky’s leefruimte beperkte zich namelijk lang tot enkel een vieze ren.
Vicky is ontzettend aanhankelijk, slim, energiek en speels. Ze houdt enorm van aandacht, maar heeft nog niet helemaal door dat je die best netjes vraagt. Aan de etiquette moet dus nog wat gewerkt worden. Je ziet Vicky heerlijk genieten van aaien, maar het allergelukkigste maak je haar door samen actief te zijn.
Samen spelen in de tuin, speuren of trainen vindt ze heerlijk.
Vicky is al goed getraind en doet het steeds beter aan de riem. Je begint wel best met wandelen in een rustige, prikkelarme omgeving. Het leven in huis kent ze nog niet, dus ze weet ook nog niet hoe ze met spullen in huis moet omgaan. Maar met het nodige geduld, leert slimme Vicky vast wel hoe ze zich als een echte dame moet gedragen.
Vicky kan niet bij andere dieren of kleine kinderen. Ze heeft ook nood aan een grote omheinde tuin om lekker in rond te rennen. 
Zoek je een energiek, aanhankelijk maatje dat je eerst vanop een veilig plekje kunt laten wennen aan het leven in huis? Dan is Vicky misschien wel de geknipte hond voor jou!
Deze hond zit in het programma van Belgian Cell Dogs, waarbij gedetineerden trainen met asielhonden om zodoende de adoptiekansen van de honden te verhogen. De honden leren commando’s en gewenst gedrag te vertonen.
Gedurende 8 weken wordt de hond getraind en in die tijd kunnen adoptanten zich heel graag aanmelden voor de hond, kunnen er adoptie gesprekken plaatsvinden én kunnen adoptanten kijken hoe er vanuit Belgian Cell Dogs wordt getraind. De adoptie vindt dan plaats na het programma van 8 weken.
'@
$ws.Range("H3").Value2 = @'
Storm
'@
$ws.Range("I3").Value2 = @'
https://dierenasielgenk.be/wp-content/uploads/2022/07/IMG_7123-1200x800.jpg
'@
$ws.Range("J3").Value2 = @'
Storm werd in beslag genomen omdat ze aan een boom werd vastgebonden. Storm heeft voordat ze in beslag werd genomen 2 andere eigenaren gehad.
Storm haar oren werken prima, wanneer zij daar zin in heeft. Wanneer ze aan het snuffelen is en ze geroepen wordt, zal ze haar gehoor even uitzetten. Zijn er lekkere snoepjes te verdienen? Kom maar op!
Deze prachtige dame is gek op eten en is daardoor ook leergierig. Er zit ook een keerzijde van haar passie voor eten: ze durft eten en objecten die zij als belangrijk ziet te beschermen.
Storm vindt het niet fijn wanneer we haar aanlijnen. We vermoeden dat haar verleden, waarbij ze aangelijnd aan een boom zat, een reden is dat ze dit gedrag inzet. Op dit moment is een halsband of harnas aandoen en wandelen niet mogelijk. Dit dient in kleine stapjes te worden opgebouwd.
Deze Siberische Husky van 2015 is actief, weet heel goed wat ze wil, en zeker ook wat ze niet wil. Ze kan niet bij andere dieren of bij kinderen. Ze heeft nood aan een fijne tuin om heerlijk in te spelen!

'@
$ws.Range("K3").Value2 = @'
Skippy
'@
$ws.Range("L3").Value2 = @'
https://dierenasielgenk.be/wp-content/uploads/2022/05/EVIOCLICK_14-10-25-1200x1800.jpg
'@
$ws.Range("M3").Value2 = @'
Skippy kwam bij ons terecht na een inbeslagname. Ze zat buiten tussen de rommel in een vieze ren.
Skippy is een heuse furie. Als ze iets wil, gaat ze er voor de honderd procent voor. Aandacht? Spelen? Ze gooit zich volledig in de strijd. In die volle overgave wordt ze soms een beetje onhandig. Zo tuimelde ze bijvoorbeeld al eens van de bank. Maar ja, hoe zou je zelf zijn als je gepassioneerd bent door speeltjes?
Het is niet omdat Skippy zo’n hevige dame is, dat ze niet kan genieten van een potje relaxen op zijn tijd. Een mand of zetel zijn héérlijke plekjes om dat te doen. Ze zal alleszins in de wolken zijn als ze binnenkort haar eigen relaxplekje heeft. Skippy kent trouwens al wat basiscommando’s zoals ZIT, AF en POOT.
Ze is een mooie Duitse Herder van 4 jaar, maar ze mag wel wat kilootjes aankomen. Ze kan niet bij andere dieren of bij kleine kinderen en ze vereist een fijne tuin om in te relaxen. 
Bied jij deze energieke herdershond een actief en liefdevol leven, dan wordt Skippy zeker en vast je beste vriendin!
'@
$ws.Range("N3").Value2 = @'
11:26:42
'@
$ws.Range("O3").Value2 = @'
11:29:44
'@
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value2 = @'
1
'@
